# Applies the "updated zeitplan and projectjournal" edit:
#  - fills in the "Effektiver Aufwand (IST)" (G) and actual start/end dates
#    (J/K) for several Zeitplan rows, letting the existing "Aufwand-
#    differenz" (H) formulas (and the G27 total) recompute automatically
#  - adds the missing H-column formula to a few rows that didn't have it yet
#  - moves the sheet selection to G20 (and scrolls so column D leads)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")
$ws.Activate()

# ---- G column (Effektiver Aufwand / IST) ---------------------------------
$gValues = @{
    9  = 4
    10 = 4
    11 = 5.5
    12 = 7
    13 = 3.5
    14 = 3
    15 = 9.5
    16 = 10.5
}
foreach ($row in $gValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $gValues[$row]
}

# ---- H column (Aufwand-differenz) -----------------------------------------
# Rows 9-15 already carry the IF(...)-formula (only their cached result was
# stale); rows 16-20 are missing the formula entirely and need it added.
foreach ($row in 16..20) {
    $ws.Cells.Item($row, 8).Formula = "=IF(G$row=`"`",`"`",SUM(`$F`$7:F$row)-SUM(`$G`$7:G$row))"
}

# ---- J / K columns (Start-/Enddatum IST) ----------------------------------
# Columns J9:J20 already use the date-formatted style, so a plain numeric
# (date-serial) assignment keeps that style untouched.
$jValues = @{
    9 = 43415; 10 = 43415; 11 = 43415; 12 = 43415; 13 = 43415
    14 = 43480; 15 = 43515; 16 = 43515; 17 = 43522
}
foreach ($row in $jValues.Keys) {
    $ws.Cells.Item($row, 10).Value = $jValues[$row]
}

# K14:K17 are still in the default "General" style in the source file, so
# pull in the date format from an already-formatted date cell (K9) before
# writing the value - this mirrors what Excel does when a user fills the
# cell down from a neighbouring, already-formatted cell.
$dateFormatSource = $ws.Range("K9")
$kValues = @{
    9 = 43473; 10 = 43473; 11 = 43473; 12 = 43473; 13 = 43473
    14 = 43480; 15 = 43550; 16 = 43529
}
foreach ($row in $kValues.Keys) {
    $target = $ws.Cells.Item($row, 11)
    $dateFormatSource.Copy()
    $target.PasteSpecial(-4122) # xlPasteFormats
    $target.Value = $kValues[$row]
}

# K17 keeps an empty value but still switches from the "General" style to
# the date style (s="15" -> s="14" in the OOXML).
$dateFormatSource.Copy()
$ws.Range("K17").PasteSpecial(-4122) # xlPasteFormats

$excel.CutCopyMode = $false

# ---- Selection / scroll position ------------------------------------------
$excel.ActiveWindow.ScrollColumn = 4   # column D leftmost
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("G20").Select()
